$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / E6: new "Appoitment" cell, styled like the rest of the row (B6) ---
$ws.Range("B6").Copy()
$ws.Range("E6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E6").Value = "Appoitment"

# --- Row 7: used to be a blank placeholder row, now a fully populated interface row ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "注册用户"
$ws.Range("C7").Value = "后台提供url，前端把注册信息传给后台，后台进行校验和数据插入，返回结果给前端"
$ws.Range("D7").Value = "注册页面"

# E7 needs the same font styling as B7/D7 (its previous style differed slightly)
$ws.Range("B6").Copy()
$ws.Range("E7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E7").Value = "客户表customer"

# Row 7 now wraps onto multiple lines because of the long C7 text
$ws.Rows("7").RowHeight = 43.2

# --- Selection moves to A8 ---
$ws.Range("A8").Select()
